$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.994.04'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.128.51'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.14'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.141.31'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.56'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.383'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.681.89'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.121.06'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.88'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.136.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000153'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '398.57'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.21'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.45'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.07'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.05'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.482'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.193'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.25%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0998'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.36%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.74'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.988'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.80'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.05'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '159.44'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.22'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.79'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.10'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.21%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.33'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.645.47'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.66'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.56%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.52'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.04'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.30'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.686'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.37%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0610'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.58%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.37'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0254'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.96'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '283.92'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.25%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0971'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.45'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.12%  '
